$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 to hold the new password value "Goku!12" (adds a new shared string)
$ws.Range("B3").Value = "Goku!12"

# Update the active cell selection to B3
$ws.Range("B3").Select()
